$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill previously-empty data cells with "no" (per commit: add "no" for empty values)
$ws.Range("L2").Value = "no"
$ws.Range("P2").Value = "no"
$ws.Range("Q2").Value = "no"
$ws.Range("AB2").Value = "no"
$ws.Range("L3").Value = "no"
$ws.Range("P3").Value = "no"
$ws.Range("Q3").Value = "no"
$ws.Range("AB3").Value = "no"
$ws.Range("L4").Value = "no"
$ws.Range("P4").Value = "no"
$ws.Range("P5").Value = "no"
$ws.Range("Q5").Value = "no"
$ws.Range("V5").Value = "no"
$ws.Range("AB5").Value = "no"
$ws.Range("L6").Value = "no"
$ws.Range("P6").Value = "no"
$ws.Range("Q6").Value = "no"
$ws.Range("AB6").Value = "no"
$ws.Range("L7").Value = "no"
$ws.Range("P7").Value = "no"
$ws.Range("L8").Value = "no"
$ws.Range("P8").Value = "no"
$ws.Range("Q8").Value = "no"
$ws.Range("AB8").Value = "no"
$ws.Range("L9").Value = "no"
$ws.Range("P9").Value = "no"
$ws.Range("Q9").Value = "no"
$ws.Range("AB9").Value = "no"
$ws.Range("L10").Value = "no"
$ws.Range("P10").Value = "no"
$ws.Range("Q10").Value = "no"
$ws.Range("AB10").Value = "no"
$ws.Range("L11").Value = "no"
$ws.Range("P11").Value = "no"
$ws.Range("V11").Value = "no"
$ws.Range("L12").Value = "no"
$ws.Range("P12").Value = "no"
$ws.Range("V12").Value = "no"
$ws.Range("L13").Value = "no"
$ws.Range("P13").Value = "no"
$ws.Range("L14").Value = "no"
$ws.Range("P14").Value = "no"
$ws.Range("Q14").Value = "no"
$ws.Range("L15").Value = "no"
$ws.Range("P15").Value = "no"
$ws.Range("V15").Value = "no"
$ws.Range("AB15").Value = "no"
$ws.Range("L16").Value = "no"
$ws.Range("V16").Value = "no"
$ws.Range("L17").Value = "no"
$ws.Range("P17").Value = "no"
$ws.Range("Q17").Value = "no"
$ws.Range("L18").Value = "no"
$ws.Range("P18").Value = "no"
$ws.Range("L19").Value = "no"
$ws.Range("P19").Value = "no"
$ws.Range("L20").Value = "no"
$ws.Range("P20").Value = "no"
$ws.Range("V20").Value = "no"
$ws.Range("P21").Value = "no"
$ws.Range("P22").Value = "no"
$ws.Range("V22").Value = "no"
$ws.Range("AB22").Value = "no"
$ws.Range("L23").Value = "no"
$ws.Range("V23").Value = "no"
$ws.Range("L24").Value = "no"
$ws.Range("V24").Value = "no"
$ws.Range("L25").Value = "no"
$ws.Range("P25").Value = "no"
$ws.Range("Q25").Value = "no"
$ws.Range("AB25").Value = "no"
$ws.Range("P26").Value = "no"
$ws.Range("Q26").Value = "no"
$ws.Range("L27").Value = "no"
$ws.Range("P27").Value = "no"
$ws.Range("Q27").Value = "no"
$ws.Range("L28").Value = "no"
$ws.Range("L29").Value = "no"
$ws.Range("L31").Value = "no"
$ws.Range("L32").Value = "no"
$ws.Range("P32").Value = "no"
$ws.Range("L33").Value = "no"
$ws.Range("P33").Value = "no"
$ws.Range("Q33").Value = "no"
$ws.Range("L34").Value = "no"
$ws.Range("P34").Value = "no"
$ws.Range("Q34").Value = "no"
$ws.Range("AB34").Value = "no"
$ws.Range("L35").Value = "no"
$ws.Range("P35").Value = "no"
$ws.Range("V35").Value = "no"
$ws.Range("AB35").Value = "no"
